$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.428.43'
$ws.Range("E2").Value = '  +2.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.080.67'
$ws.Range("E3").Value = '  +2.53%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.49'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.36'
$ws.Range("E7").Value = '  +1.24%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +1.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0833'
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.386.53'
$ws.Range("E12").Value = '  +2.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.98'
$ws.Range("E13").Value = '  +4.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.34'
$ws.Range("E14").Value = '  +6.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.781'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.45'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.070.61'
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.416.12'
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("E19").Value = '  +2.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.01'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0832'
$ws.Range("E21").Value = '  +1.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.16'
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.85'
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.44'
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.136'
$ws.Range("E28").Value = '  +6.24%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.38'
$ws.Range("E29").Value = '  +9.32%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.07'
$ws.Range("E30").Value = '  +1.92%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("E32").Value = '  +4.47%  '
$ws.Range("E33").Value = '  +7.45%  '
$ws.Range("E34").Value = '  +2.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0608'
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.36'
$ws.Range("E36").Value = '  +1.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.31'
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("E38").Value = '  +3.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.30'
$ws.Range("E40").Value = '  +2.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.538.24'
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.16'
$ws.Range("E42").Value = '  +3.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0221'
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("E44").Value = '  -0.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0922'
$ws.Range("E45").Value = '  +1.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.65'
$ws.Range("E46").Value = '  +8.69%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.12'
$ws.Range("E47").Value = '  -1.33%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.11'
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("E49").Value = '  +2.67%  '
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.274.97'
$ws.Range("E51").Value = '  +2.20%  '
